$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" metadata property value (row 8, column B)
$ws.Range("B8").Value = "2023-04-12T13:10:15+00:00"

# Fill in the "Case Sensitive" metadata property value (row 14, column B)
# with the text "true". A leading apostrophe forces Excel to store it as
# text instead of auto-converting it to a Boolean, and PasteSpecial
# (formats only) from the neighbouring row restores the normal cell style
# that the apostrophe-prefix entry mode would otherwise override.
$cell = $ws.Range("B14")
$cell.Value = "'true"

$styleSource = $ws.Range("B15")
$styleSource.Copy()
$cell.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
